# Build site at 2022-01-09 00:29:46 UTC
# Adds a new "Docentes responsáveis:" / "8855158 - Morun Bernardino Neto"
# entry right after the Objectives rows (old rows 12-23 shift down to 14-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the old row 12 ("Programa resumido:").
# Inserting whole rows shifts everything below down and carries the row-style
# (s="1"/"2"/"3") from the row immediately above, matching the target layout.
$ws.Rows("12:13").Insert()

# Populate the newly inserted rows (the remaining cells in these two rows
# - B12/C12 and A13 - stay blank, matching the target layout).
$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Range("B13").Value = "8855158 - Morun Bernardino Neto"
$ws.Range("C13").Value = "8855158 - Morun Bernardino Neto"
